$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "p24_1" column (Z) entirely - header and all data
$ws.Range("Z:Z").Delete()

# "nota_iniciativa" column (Y) data rows: normalize every value to 0
$lastRow = $ws.Range("A1").SpecialCells(11).Row
$ws.Range("Y2:Y" + $lastRow).Value = 0
